# Update for Report node structure
$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "PositiveExtra"
$ws2 = $wb.Worksheets.Item(2)   # "ProductList"

# --- PositiveExtra sheet ---
$ws1.Activate() | Out-Null

# B2:B15 values 13.5 -> 10.5
$ws1.Range("B2:B15").Value = 10.5

# Move the selection/view on PositiveExtra from A2:A15 to B2:B15
$ws1.Range("B2:B15").Select() | Out-Null

# --- ProductList sheet becomes the active / selected tab, keeping its C14 selection ---
$ws2.Activate() | Out-Null
$ws2.Range("C14").Select() | Out-Null
